# Uplift to the new version of pyxform/cht-conf
#
# 1) "survey" sheet: remove the orphan "NO_LABEL" note that used to sit in C3
#    (D3 keeps its "field-list" value).
# 2) "settings" sheet: drop the "form_id" column (column B) entirely, which
#    shifts version/style/namespaces (and their sample values) one column to
#    the left. The column header comments need to carry their text over to
#    the new column that ends up holding that topic.
# 3) Tidy up the "settings" sheet selection to point at the new layout.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- 1) survey sheet: remove the stray NO_LABEL value from C3 ---------------
$survey.Range("C3").Clear()

# --- 2) settings sheet: re-home the header comments before the column shift -
# Grab the original text from the columns that will end up one slot to the
# left once form_id (column B) disappears, then write it onto the surviving
# header cells. The E1 (namespaces) comment becomes D1; there is no longer an
# E1 header, so that original comment is simply dropped.
$versionComment = $settings.Range("C1").Comment.Text()
$styleComment = $settings.Range("D1").Comment.Text()
$namespacesComment = $settings.Range("E1").Comment.Text()

$settings.Range("E1").Comment.Delete()
$settings.Range("D1").Comment.Text($namespacesComment)
$settings.Range("C1").Comment.Text($styleComment)
$settings.Range("B1").Comment.Text($versionComment)

# --- now remove the form_id column (B) -- everything right of it shifts left
$settings.Columns.Item(2).Delete()

# --- 3) refresh selections to match the new layout ---------------------------
$settings.Range("B1").Select()
$survey.Range("A2").Select()
